# chi_square_results_v2.xlsx was reopened/resaved in a newer Excel build.
# The underlying data/formulas are unchanged (feature ranking ensembled with
# top 10 feature ranks); what actually changed on the worksheet itself is:
#   - column A / column B given explicit (best-fit) widths
#   - the live selection left on K34 after the user finished reviewing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-fit column widths (A: feature name column, B: chi-square value column)
$ws.Columns.Item(1).ColumnWidth = 22.75
$ws.Columns.Item(2).ColumnWidth = 13.42

# Leave the selection where the author last left it
$ws.Range("K34").Select()
